{"js": "// \"Better handle copying paragraph styles\"\n//\n// When a \"join\" paragraph (style MSC_Join, used to glue two passage\n// chunks together) is expanded into the \"[...]\" placeholder paragraph\n// plus its trailing blank paragraph, those two generated paragraphs\n// must carry the same MSC_Join paragraph style as the paragraph they\n// were copied from, instead of being left with the default style.\n//\n// Find every \"MSC_Join\" paragraph that is immediately followed by an\n// empty \"[...]\" placeholder paragraph and then a blank paragraph, and\n// apply the MSC_Join style to those two following paragraphs as well.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst JOIN_STYLE = \"MSC_Join\";\n\nfor (let i = 0; i < items.length; i++) {\n  const joinPara = items[i];\n  if (joinPara.style !== JOIN_STYLE || joinPara.text !== \"\") {\n    continue;\n  }\n\n  const placeholderPara = items[i + 1];\n  const blankPara = items[i + 2];\n  if (!placeholderPara || !blankPara) {\n    continue;\n  }\n  if (placeholderPara.text !== \"[...]\" || blankPara.text !== \"\") {\n    continue;\n  }\n\n  placeholderPara.style = JOIN_STYLE;\n  blankPara.style = JOIN_STYLE;\n}\n\nawait context.sync();\n", "ps1": "# \"Better handle copying paragraph styles\"\n#\n# When a \"join\" paragraph (style MSC_Join, used to glue two passage\n# chunks together) is expanded into the \"[...]\" placeholder paragraph\n# plus its trailing blank paragraph, those two generated paragraphs\n# must carry the same MSC_Join paragraph style as the paragraph they\n# were copied from, instead of being left with the default style.\n#\n# Find every \"MSC_Join\" paragraph that is immediately followed by an\n# empty \"[...]\" placeholder paragraph and then a blank paragraph, and\n# apply the MSC_Join style to those two following paragraphs as well.\n\n$d = $word.ActiveDocument\n$JoinStyle = \"MSC_Join\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count - 2; $i++) {\n    $joinPara = $d.Paragraphs.Item($i)\n    if ($joinPara.Style.NameLocal -ne $JoinStyle) {\n        continue\n    }\n    if ($joinPara.Range.Text.Trim() -ne \"\") {\n        continue\n    }\n\n    $placeholderPara = $d.Paragraphs.Item($i + 1)\n    $blankPara = $d.Paragraphs.Item($i + 2)\n\n    if ($placeholderPara.Range.Text.Trim() -ne \"[...]\") {\n        continue\n    }\n    if ($blankPara.Range.Text.Trim() -ne \"\") {\n        continue\n    }\n\n    $placeholderPara.Style = $JoinStyle\n    $blankPara.Style = $JoinStyle\n}\n"}
